# "Generate Report for Archive" — refresh the localization-status report:
#   1. Every cell that used to read "Ready for handoff" now reads "In Translation".
#   2. The "zh-cn"/"de-de" status columns on the Overview sheet (E:F) and the
#      "Status" column (C) on each per-language sheet are narrowed to match
#      the regenerated report's column sizing.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update status text wherever it said "Ready for handoff" ---
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# --- 2. Narrow the status columns ---
# ColumnWidth is expressed in characters of the Normal style font and Excel
# snaps it to the nearest whole on-screen pixel, so we pick the character
# width whose rounded column matches the report's regenerated column size.
$overview.Columns.Item(5).ColumnWidth = 13.2
$overview.Columns.Item(6).ColumnWidth = 13.2
$zhcn.Columns.Item(3).ColumnWidth     = 13.2
$dede.Columns.Item(3).ColumnWidth     = 13.2
